# "Generate Report for Handback" - update the localization-status report
# after a handback (target files generated, handback timestamps recorded).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1) Status text: every "Ready for handoff" cell becomes "Handed back: in
#    sync with en-US" (Overview!E2:F2/E3:F3, zh-cn!C2:C3, de-de!C2:C3).
# ---------------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
}

# ---------------------------------------------------------------------------
# 2) Latest Handback DateTime: previously-unset "0001-01-01 00:00:00" is
#    replaced everywhere with the zh-cn handback stamp, then de-de's own
#    K2/K3 are overwritten with its own (slightly later) handback stamp.
# ---------------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("0001-01-01 00:00:00", "2016-08-23 00:45:24")
}
$wsDeDe.Range("K2").Value = "2016-08-23 00:45:31"
$wsDeDe.Range("K3").Value = "2016-08-23 00:45:31"

# ---------------------------------------------------------------------------
# 3) Latest Target File / Latest Handback File columns (I, J) now get
#    populated for both language sheets, rows 2 and 3. Column I becomes a
#    hyperlink to the same markdown source as column A; column J holds the
#    generated handback xliff file name.
# ---------------------------------------------------------------------------
$zhCn36038Target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c664e18cfcefa78c8e3f5de7989a7b424ee7775b/e2e/36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$zhCnFbebf6Target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c664e18cfcefa78c8e3f5de7989a7b424ee7775b/e2e/fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"

# zh-cn
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhCn36038Target, "", "", "36038c19-52ce-4a1b-8036-de19daaeacb8.md")
$wsZhCn.Range("J2").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.4ad38bf24ece5dede3f5c797292c36779823b837.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhCnFbebf6Target, "", "", "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md")
$wsZhCn.Range("J3").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.3180788ccb1d5c858ef0ef8e59d53f0fa210ab48.zh-cn.xlf"

# de-de
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $zhCn36038Target, "", "", "36038c19-52ce-4a1b-8036-de19daaeacb8.md")
$wsDeDe.Range("J2").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.4ad38bf24ece5dede3f5c797292c36779823b837.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $zhCnFbebf6Target, "", "", "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md")
$wsDeDe.Range("J3").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.3180788ccb1d5c858ef0ef8e59d53f0fa210ab48.de-de.xlf"

# ---------------------------------------------------------------------------
# 4) Column widths grow to fit the newly-populated / re-worded columns.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.15   # E
$wsOverview.Columns.Item(6).ColumnWidth = 29.15   # F

$wsZhCn.Columns.Item(3).ColumnWidth = 29.15        # C  Status
$wsZhCn.Columns.Item(9).ColumnWidth = 39.15        # I  Latest Target File
$wsZhCn.Columns.Item(10).ColumnWidth = 39.15       # J  Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth = 29.15        # C  Status
$wsDeDe.Columns.Item(9).ColumnWidth = 39.15        # I  Latest Target File
$wsDeDe.Columns.Item(10).ColumnWidth = 39.15       # J  Latest Handback File
